# ustore_inventory.xlsx — "added building list, more cart functionality"
#
# Content change: the "School Supplies" category (rows 129-145, column C)
# is renamed to "School". (TigerTask app's shared-string table ends up
# dropping the old "School Supplies" string and appending a new "School"
# string at the end, which is what Excel does automatically when a shared
# string's last reference is replaced with different text — the item
# names in column A are untouched.)
#
# View change: the active selection moves to C137 (and the window was
# scrolled so A119 is the top-left visible cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "School Supplies" category to "School" for the school-supplies
# block of rows (129-145).
$ws.Range("C129:C145").Value = "School"

# Update the saved selection/view state to match the edited workbook.
$win = $wb.Windows.Item(1)
$ws.Range("C137").Select()
$win.ScrollRow = 119
$win.ScrollColumn = 1
